# Trading update: 2026-02-18 10:17:53
#
# This script applies the latest trading-bot snapshot update to the
# live_trading_results workbook:
#   - Summary: Current Capital rounds down from 1500.25 to 1500
#   - Strategy Status: MarketMaking strategy capital/PnL% reset to 100 / 0
#   - All Trades: header columns re-ordered to match the per-strategy sheet
#     layout, older rows lose their "latest snapshot" enrichment columns,
#     two new trades (#6, #7) are appended
#   - MarketMaking: the per-strategy sheet now tracks trade #7 (the latest
#     open MarketMaking trade) instead of trade #4, and the stale trade #5
#     row is removed

$wb = $excel.ActiveWorkbook

# Helper: write a text value into a cell without Excel's automatic
# date/number detection turning it into a different type. A leading
# apostrophe forces "treat as text"; ClearFormats() then strips the
# resulting quote-prefix style so the cell is left without any explicit
# style, matching the rest of the sheet.
function Set-TextValue {
    param($range, [string]$text)
    $range.Value = "'" + $text
    $range.ClearFormats()
}

# Helper: turn a cell into an empty (but present) text cell, matching the
# workbook's convention of blank placeholder cells (e.g. N3/O3/P3/Q3 in
# "All Trades"), rather than removing the cell outright.
function Set-BlankMarker {
    param($range)
    $range.Value = "'"
    $range.ClearFormats()
}

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1500

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C6").Value = 100
$status.Range("F6").Value = 0

# ---------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------
$all = $wb.Worksheets.Item("All Trades")

# Header columns L:Q re-ordered to match the MarketMaking sheet layout
$all.Range("L1").Value = "Entry Slippage (bps)"
$all.Range("M1").Value = "Exit Slippage (bps)"
$all.Range("N1").Value = "Confidence"
$all.Range("O1").Value = "Entry Reason"
$all.Range("P1").Value = "Exit Reason"
$all.Range("Q1").Value = "Duration (min)"

# Row 3 (trade #2) loses its "latest snapshot" enrichment
Set-BlankMarker $all.Range("K3")
Set-BlankMarker $all.Range("L3")
Set-BlankMarker $all.Range("M3")

# Row 5 (trade #4) - exit price now recorded as 0, loses enrichment columns
$all.Range("G5").Value = 0
Set-BlankMarker $all.Range("K5")
Set-BlankMarker $all.Range("M5")
Set-BlankMarker $all.Range("N5")
Set-BlankMarker $all.Range("O5")
Set-BlankMarker $all.Range("P5")
Set-BlankMarker $all.Range("Q5")

# Row 6 (trade #5) - exit price now recorded as 0, loses enrichment columns
$all.Range("G6").Value = 0
Set-BlankMarker $all.Range("K6")
Set-BlankMarker $all.Range("L6")
Set-BlankMarker $all.Range("M6")
Set-BlankMarker $all.Range("N6")
Set-BlankMarker $all.Range("O6")
Set-BlankMarker $all.Range("P6")
Set-BlankMarker $all.Range("Q6")

# Row 7 (NEW) - trade #6
$all.Range("A7").Value = 6
Set-TextValue $all.Range("B7") "2026-02-18"
Set-TextValue $all.Range("C7") "10:16:20"
$all.Range("D7").Value = "MarketMaking"
$all.Range("E7").Value = "DOWN"
$all.Range("F7").Value = 0.43
$all.Range("G7").Value = 0
$all.Range("H7").Value = "OPEN"
$all.Range("I7").Value = 0
$all.Range("J7").Value = 0
Set-BlankMarker $all.Range("K7")
Set-BlankMarker $all.Range("L7")
Set-BlankMarker $all.Range("M7")
Set-BlankMarker $all.Range("N7")
Set-BlankMarker $all.Range("O7")
Set-BlankMarker $all.Range("P7")
Set-BlankMarker $all.Range("Q7")

# Row 8 (NEW) - trade #7, the latest trade, carries the enrichment columns
$all.Range("A8").Value = 7
Set-TextValue $all.Range("B8") "2026-02-18"
Set-TextValue $all.Range("C8") "10:17:35"
$all.Range("D8").Value = "MarketMaking"
$all.Range("E8").Value = "UP"
$all.Range("F8").Value = 0.67
Set-BlankMarker $all.Range("G8")
$all.Range("H8").Value = "OPEN"
$all.Range("I8").Value = 0
$all.Range("J8").Value = 0
$all.Range("K8").Value = 100
$all.Range("L8").Value = 0
$all.Range("M8").Value = 0
$all.Range("N8").Value = 0.6
$all.Range("O8").Value = "Normal spread capture: 202 bps"
Set-BlankMarker $all.Range("P8")
$all.Range("Q8").Value = 0

# ---------------------------------------------------------------
# MarketMaking sheet - now tracks trade #7 instead of trade #4,
# and the stale trade #5 row is dropped.
# ---------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("A2").Value = 7
Set-TextValue $mm.Range("C2") "10:17:35"
$mm.Range("F2").Value = 0.67
$mm.Rows.Item(3).Delete()
